# Updates the "cryptos" price/volume table to reflect the latest scrape.
# Values are written so that they remain plain text cells (matching the
# original inlineStr cells), even for strings that look numeric/date-like
# (e.g. "1.00", "4.69"), and without leaving a residual explicit cell
# style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $value) {
    $range = $ws.Range($cellRef)
    # Force a text number format first so Excel does not reinterpret the
    # assigned string as a number/date (which would e.g. turn "1.00" into 1
    # or "48.377.79" into some other representation).
    $range.NumberFormat = "@"
    $range.Value = $value
    # Restore the default "Normal" style so we don't leave a stray style
    # index on the cell compared to the original (unstyled) cells.
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws "D2" "48.377.79"
Set-TextValue $ws "E2" "  +1.52%  "

# Row 3 - Ethereum
Set-TextValue $ws "D3" "2.511.50"
Set-TextValue $ws "E3" "  +0.68%  "

# Row 4 - TetherUSD
Set-TextValue $ws "D4" "1.00"
Set-TextValue $ws "E4" "  +0.06%  "

# Row 5 - BNB
Set-TextValue $ws "D5" "321.56"
Set-TextValue $ws "E5" "  -0.04%  "

# Row 6 - Solana
Set-TextValue $ws "D6" "108.61"
Set-TextValue $ws "E6" "  -0.39%  "

# Row 7 - XRP
Set-TextValue $ws "D7" "0.529"
Set-TextValue $ws "E7" "  +1.10%  "

# Row 8 - USDC
Set-TextValue $ws "D8" "1.00"
Set-TextValue $ws "E8" "  +0.03%  "

# Row 9 - Cardano
Set-TextValue $ws "E9" "  +0.06%  "

# Row 10 - Avalanche
Set-TextValue $ws "D10" "40.02"
Set-TextValue $ws "E10" "  +1.67%  "

# Row 11 - Chainlink
Set-TextValue $ws "D11" "20.33"
Set-TextValue $ws "E11" "  +8.93%  "

# Row 12 - Dogecoin
Set-TextValue $ws "E12" "  +1.08%  "

# Row 14 - Polkadot
Set-TextValue $ws "E14" "  -0.42%  "

# Row 15 - Wrapped liquid staked Ether 2.0
Set-TextValue $ws "D15" "2.903.45"
Set-TextValue $ws "E15" "  +0.67%  "

# Row 16 - Wrapped Ether
Set-TextValue $ws "D16" "2.510.92"
Set-TextValue $ws "E16" "  +0.71%  "

# Row 17 - Polygon
Set-TextValue $ws "E17" "  -0.27%  "

# Row 18 - Wrapped BTC
Set-TextValue $ws "D18" "48.222.21"
Set-TextValue $ws "E18" "  +1.53%  "

# Row 19 - Internet Computer (DFINITY)
Set-TextValue $ws "D19" "13.15"
Set-TextValue $ws "E19" "  -1.33%  "

# Row 20 - Uniswap
Set-TextValue $ws "D20" "6.79"
Set-TextValue $ws "E20" "  +2.38%  "

# Row 21 - Shiba Inu
Set-TextValue $ws "E21" "  +0.12%  "

# Row 22 - ImmutableX
Set-TextValue $ws "E22" "  +0.42%  "

# Row 23 - Bitcoin Cash
Set-TextValue $ws "D23" "279.50"
Set-TextValue $ws "E23" "  +13.09%  "

# Row 24 - Litecoin
Set-TextValue $ws "D24" "72.25"
Set-TextValue $ws "E24" "  +2.27%  "

# Row 25 - PancakeSwap
Set-TextValue $ws "D25" "2.55"
Set-TextValue $ws "E25" "  -0.08%  "

# Row 26 - Dai
Set-TextValue $ws "E26" "  -0.05%  "

# Row 27 - Ethereum Classic
Set-TextValue $ws "D27" "25.79"
Set-TextValue $ws "E27" "  +0.02%  "

# Row 28 - Toncoin
Set-TextValue $ws "E28" "  -4.34%  "

# Row 29 - Cosmos
Set-TextValue $ws "E29" "  -1.96%  "

# Row 30 - Kaspa
Set-TextValue $ws "D30" "0.141"
Set-TextValue $ws "E30" "  +1.44%  "

# Row 31 - Injective Protocol
Set-TextValue $ws "D31" "35.37"
Set-TextValue $ws "E31" "  +1.96%  "

# Row 32 - OKB
Set-TextValue $ws "D32" "49.60"
Set-TextValue $ws "E32" "  -0.67%  "

# Row 33 - Celestia
Set-TextValue $ws "D33" "19.77"
Set-TextValue $ws "E33" "  -2.11%  "

# Row 34 - Filecoin
Set-TextValue $ws "D34" "5.36"
Set-TextValue $ws "E34" "  +0.44%  "

# Row 35 - First Digital USD
Set-TextValue $ws "E35" "  -0.03%  "

# Row 36 - Hedera
Set-TextValue $ws "D36" "0.0784"
Set-TextValue $ws "E36" "  -0.37%  "

# Row 37 / 38 - these two rows swapped places (RenderToken <-> ARBITRUM)
Set-TextValue $ws "B37" "RenderToken"
Set-TextValue $ws "C37" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D37" "4.69"
Set-TextValue $ws "E37" "  -0.80%  "

Set-TextValue $ws "B38" "ARBITRUM"
Set-TextValue $ws "C38" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws "D38" "1.96"
Set-TextValue $ws "E38" "  +0.00%  "

# Row 39 - LidoDAOToken
Set-TextValue $ws "E39" "  -1.01%  "

# Row 40 - Stellar
Set-TextValue $ws "E40" "  -0.18%  "

# Row 41 - Monero
Set-TextValue $ws "D41" "121.95"
Set-TextValue $ws "E41" "  +2.05%  "

# Row 42 - WEMIX Token
Set-TextValue $ws "E42" "  +0.18%  "

# Row 43 - EnergySwap
Set-TextValue $ws "D43" "21.62"
Set-TextValue $ws "E43" "  -2.45%  "

# Row 44 - VeChain
Set-TextValue $ws "E44" "  +1.79%  "

# Row 45 - Maker
Set-TextValue $ws "D45" "2.020.40"
Set-TextValue $ws "E45" "  +1.36%  "

# Row 46 - NEAR Protocol
Set-TextValue $ws "D46" "3.21"
Set-TextValue $ws "E46" "  +6.06%  "

# Row 47 - Stacks
Set-TextValue $ws "E47" "  +2.74%  "

# Row 48 - ApeX Protocol
Set-TextValue $ws "D48" "2.00"
Set-TextValue $ws "E48" "  -2.30%  "

# Row 50 - THORChain
Set-TextValue $ws "D50" "5.20"
Set-TextValue $ws "E50" "  -0.51%  "

# Row 51 - Bitcoin SV
Set-TextValue $ws "D51" "80.75"
Set-TextValue $ws "E51" "  +3.99%  "
